# إضافة حدث جديد في Card24
# Adds a new service-log row (row 16) to the "Card24" sheet and backfills the
# previously-blank cells on row 15 with the literal text "nan" (matching the
# sheet's existing convention for empty data cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Card24")

# --- Row 15: fill the previously-empty cells with the literal text "nan" ---
$ws.Range("B15").Value = "nan"
$ws.Range("C15").Value = "nan"
$ws.Range("E15").Value = "nan"
$ws.Range("F15").Value = "nan"
$ws.Range("G15").Value = "nan"
$ws.Range("H15").Value = "nan"
$ws.Range("I15").Value = "nan"
$ws.Range("J15").Value = "nan"
$ws.Range("K15").Value = "nan"
$ws.Range("P15").Value = "nan"

# --- Row 16: brand-new service-log entry for card 24 ---
# Numeric-looking values are entered with a leading apostrophe so Excel keeps
# them as plain text (matching the rest of the sheet, which stores every
# value - even numbers and dates - as text).
$ws.Range("A16").Value = "'24"
$ws.Range("D16").Value = "'968"
$ws.Range("L16").Value = "'11/12/2025"
$ws.Range("M16").Value = "زياره توكيل"
$ws.Range("N16").Value = "تم  سن دوغر وسلندر وفلاتس وعيار ماكينه"
$ws.Range("O16").Value = "خبير ارول"
